# Update automàtic: dades i banners [2026-02-08 18:50]
# Refreshes the per-station extraction timestamps and the latest observed
# measurements (snow depth, precipitation, pressure, humidity, wind gust,
# temperatures) pulled from meteo.cat for 2026-02-08.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell values keyed by A1 reference.
$updates = @{
    'E2' = '2026-02-08 18:48:21'
    'I2' = '4.0 mm'
    'E3' = '2026-02-08 18:48:23'
    'G3' = '173 cm'
    'I3' = '2.1 mm'
    'E4' = '2026-02-08 18:48:26'
    'I4' = '1.6 mm'
    'J4' = '1001.8 hPa'
    'E5' = '2026-02-08 18:48:28'
    'G5' = '118 cm'
    'E6' = '2026-02-08 18:48:30'
    'J6' = '1001.6 hPa'
    'E7' = '2026-02-08 18:48:33'
    'J7' = '1002.0 hPa'
    'E8' = '2026-02-08 18:48:35'
    'J8' = '1001.9 hPa'
    'O8' = '9.5 °C'
    'E9' = '2026-02-08 18:48:38'
    'E10' = '2026-02-08 18:48:40'
    'E11' = '2026-02-08 18:48:43'
    'O11' = '4.8 °C'
    'E12' = '2026-02-08 18:48:45'
    'E13' = '2026-02-08 18:48:47'
    'J13' = '1003.4 hPa'
    'O13' = '3.5 °C'
    'E14' = '2026-02-08 18:48:50'
    'E15' = '2026-02-08 18:48:52'
    'E16' = '2026-02-08 18:48:54'
    'I16' = '3.0 mm'
    'E17' = '2026-02-08 18:48:57'
    'O17' = '0.0 °C'
    'E18' = '2026-02-08 18:48:59'
    'J18' = '1002.0 hPa'
    'E19' = '2026-02-08 18:49:02'
    'I19' = '10.1 mm'
    'E20' = '2026-02-08 18:49:04'
    'I20' = '7.5 mm'
    'L20' = '50.8 km/h - 338º 18:28 TU'
    'O20' = '-4.3 °C'
    'E21' = '2026-02-08 18:49:07'
    'H21' = '80%'
    'J21' = '1002.9 hPa'
    'E22' = '2026-02-08 18:49:09'
    'E23' = '2026-02-08 18:49:12'
    'I23' = '4.3 mm'
    'E24' = '2026-02-08 18:49:14'
    'H24' = '82%'
    'J24' = '1003.3 hPa'
    'E25' = '2026-02-08 18:49:17'
    'I25' = '0.6 mm'
    'E26' = '2026-02-08 18:49:19'
    'H26' = '67%'
    'J26' = '1001.0 hPa'
    'O26' = '3.8 °C'
    'E27' = '2026-02-08 18:49:22'
    'E28' = '2026-02-08 18:49:24'
    'J28' = '1001.6 hPa'
    'E29' = '2026-02-08 18:49:26'
    'E30' = '2026-02-08 18:49:29'
    'J30' = '1002.0 hPa'
    'E31' = '2026-02-08 18:49:31'
    'H31' = '74%'
    'J31' = '1001.1 hPa'
    'N31' = '8.5 °C 18:27 TU'
    'E32' = '2026-02-08 18:49:34'
    'E33' = '2026-02-08 18:49:37'
    'J33' = '1003.1 hPa'
    'O33' = '2.8 °C'
    'E34' = '2026-02-08 18:49:39'
    'E35' = '2026-02-08 18:49:41'
    'J35' = '1004.0 hPa'
    'E36' = '2026-02-08 18:49:44'
    'J36' = '1002.1 hPa'
    'O36' = '11.2 °C'
    'E37' = '2026-02-08 18:49:46'
    'J37' = '1002.9 hPa'
    'E38' = '2026-02-08 18:49:49'
    'I38' = '0.6 mm'
    'E39' = '2026-02-08 18:49:51'
    'E40' = '2026-02-08 18:49:53'
    'J40' = '1003.4 hPa'
    'E41' = '2026-02-08 18:49:55'
    'H41' = '71%'
    'J41' = '1002.1 hPa'
    'O41' = '12.0 °C'
    'E42' = '2026-02-08 18:49:58'
    'E43' = '2026-02-08 18:50:00'
    'O43' = '6.9 °C'
    'E44' = '2026-02-08 18:50:03'
    'I44' = '2.1 mm'
    'E45' = '2026-02-08 18:50:05'
    'H45' = '76%'
    'J45' = '1004.2 hPa'
    'E46' = '2026-02-08 18:50:07'
    'H46' = '74%'
    'J46' = '1003.8 hPa'
    'O46' = '9.5 °C'
}

# Cells whose new text looks like "NN%" need to be force-formatted as Text
# first, otherwise Excel's input parser silently turns "80%" into the
# number 0.8 displayed with a Percentage format instead of keeping the
# literal string "80%" that the source data uses.
$percentCells = @('H21', 'H24', 'H26', 'H31', 'H41', 'H45', 'H46')
foreach ($pc in $percentCells) {
    $ws.Range($pc).NumberFormat = "@"
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
